$wb = $excel.ActiveWorkbook

# --- ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 121.87
$ws.Range("I15").Value = 121.87
$ws.Range("K15").Value = 365.61
$ws.Range("M15").Value = -196.61
$ws.Range("H98").Value = 1952
$ws.Range("I98").Value = 1444
$ws.Range("K98").Value = 1444
$ws.Range("M98").Value = 54
$ws.Range("H107").Value = 741.5
$ws.Range("I107").Value = 569.5833
$ws.Range("J107").Value = 947.8
$ws.Range("K107").Value = 569.5833
$ws.Range("L107").Value = 947.8
$ws.Range("M107").Value = 1350.4167
$ws.Range("N107").Value = -4787.8
$ws.Range("H122").Value = 1952
$ws.Range("I122").Value = 1444
$ws.Range("K122").Value = 4332
$ws.Range("M122").Value = -1882
$ws.Range("H138").Value = 3708028.8
$ws.Range("I138").Value = 3356.0715
$ws.Range("J138").Value = 5004664
$ws.Range("K138").Value = 10068.2145
$ws.Range("L138").Value = 15013992
$ws.Range("M138").Value = -4928.2145
$ws.Range("N138").Value = -15024272

# --- ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 6854490
$ws.Range("I32").Value = 53611.312
$ws.Range("J32").Value = 19656144
$ws.Range("K32").Value = 53611.312
$ws.Range("L32").Value = 19656144
$ws.Range("M32").Value = -53324.312
$ws.Range("N32").Value = -19656718
$ws.Range("H61").Value = 30366148
$ws.Range("I61").Value = 40041704
$ws.Range("K61").Value = 40041704
$ws.Range("M61").Value = -40041492
$ws.Range("H74").Value = 10481729
$ws.Range("I74").Value = 16718571
$ws.Range("J74").Value = 86991.836
$ws.Range("K74").Value = 16718571
$ws.Range("L74").Value = 86991.836
$ws.Range("M74").Value = -16717697
$ws.Range("N74").Value = -88739.836
$ws.Range("H77").Value = 10481729
$ws.Range("I77").Value = 16718571
$ws.Range("J77").Value = 86991.836
$ws.Range("K77").Value = 83592855
$ws.Range("L77").Value = 434959.18
$ws.Range("M77").Value = -83588487
$ws.Range("N77").Value = -443695.18
$ws.Range("H122").Value = 3586388.5
$ws.Range("I122").Value = 1664.7894
$ws.Range("J122").Value = 9262201
$ws.Range("K122").Value = 4994.3682
$ws.Range("L122").Value = 27786603
$ws.Range("M122").Value = -2544.3682
$ws.Range("N122").Value = -27791503
$ws.Range("H132").Value = 7971238
$ws.Range("I132").Value = 10225894
$ws.Range("J132").Value = 79943.92999999999
$ws.Range("K132").Value = 30677682
$ws.Range("L132").Value = 239831.79
$ws.Range("M132").Value = -30675152
$ws.Range("N132").Value = -244891.79
$ws.Range("H136").Value = 30366148
$ws.Range("I136").Value = 40041704
$ws.Range("K136").Value = 120125112
$ws.Range("M136").Value = -120122562

# --- CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 54668.15
$ws.Range("I31").Value = 33553.938
$ws.Range("J31").Value = 139125
$ws.Range("K31").Value = 33553.938
$ws.Range("L31").Value = 139125
$ws.Range("M31").Value = -33258.938
$ws.Range("N31").Value = -139715
$ws.Range("H34").Value = 54668.15
$ws.Range("I34").Value = 33553.938
$ws.Range("J34").Value = 139125
$ws.Range("K34").Value = 33553.938
$ws.Range("L34").Value = 139125
$ws.Range("M34").Value = -33351.938
$ws.Range("N34").Value = -139529
$ws.Range("H58").Value = 47620456
$ws.Range("I58").Value = 52632956
$ws.Range("J58").Value = 1700
$ws.Range("K58").Value = 52632956
$ws.Range("L58").Value = 1700
$ws.Range("M58").Value = -52632753
$ws.Range("N58").Value = -2106
$ws.Range("H122").Value = 2228
$ws.Range("I122").Value = 1456
$ws.Range("J122").Value = 3000
$ws.Range("K122").Value = 4368
$ws.Range("L122").Value = 9000
$ws.Range("M122").Value = -1918
$ws.Range("N122").Value = -13900
$ws.Range("H132").Value = 50475.76
$ws.Range("I132").Value = 3126.6365
$ws.Range("J132").Value = 102559.8
$ws.Range("K132").Value = 9379.9095
$ws.Range("L132").Value = 307679.4
$ws.Range("M132").Value = -6849.9095
$ws.Range("N132").Value = -312739.4
$ws.Range("H136").Value = 47620456
$ws.Range("I136").Value = 52632956
$ws.Range("J136").Value = 1700
$ws.Range("K136").Value = 157898868
$ws.Range("L136").Value = 5100
$ws.Range("M136").Value = -157896318
$ws.Range("N136").Value = -10200

# --- CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H130").Value = 2683.7273
$ws.Range("J130").Value = 3323.75
$ws.Range("L130").Value = 9971.25
$ws.Range("N130").Value = -20011.25
$ws.Range("H131").Value = 1042.7126
$ws.Range("J131").Value = 1082.9136
$ws.Range("L131").Value = 3248.7408
$ws.Range("N131").Value = -13328.7408
$ws.Range("H132").Value = 2273.0952
$ws.Range("I132").Value = 1554.9333
$ws.Range("J132").Value = 2672.074
$ws.Range("K132").Value = 13994.3997
$ws.Range("L132").Value = 24048.666
$ws.Range("M132").Value = -11464.3997
$ws.Range("N132").Value = -29108.666
$ws.Range("H137").Value = 3072.65
$ws.Range("I137").Value = 1173.3334
$ws.Range("J137").Value = 3407.8235
$ws.Range("K137").Value = 3520.0002
$ws.Range("L137").Value = 10223.4705
$ws.Range("M137").Value = 1579.9998
$ws.Range("N137").Value = -20423.4705

# --- LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1089.3334
$ws.Range("I22").Value = 850.5
$ws.Range("J22").Value = 1280.4
$ws.Range("K22").Value = 850.5
$ws.Range("L22").Value = 1280.4
$ws.Range("M22").Value = -555.5
$ws.Range("N22").Value = -1870.4
$ws.Range("H27").Value = 1089.3334
$ws.Range("I27").Value = 850.5
$ws.Range("J27").Value = 1280.4
$ws.Range("K27").Value = 850.5
$ws.Range("L27").Value = 1280.4
$ws.Range("M27").Value = -743.5
$ws.Range("N27").Value = -1494.4
$ws.Range("H61").Value = 2827.6365
$ws.Range("I61").Value = 2806.125
$ws.Range("K61").Value = 2806.125
$ws.Range("M61").Value = -2604.125
$ws.Range("H82").Value = 1326.5
$ws.Range("I82").Value = 1059.25
$ws.Range("J82").Value = 1861
$ws.Range("K82").Value = 1059.25
$ws.Range("L82").Value = 1861
$ws.Range("M82").Value = -698.25
$ws.Range("N82").Value = -2583
$ws.Range("H85").Value = 1326.5
$ws.Range("I85").Value = 1059.25
$ws.Range("J85").Value = 1861
$ws.Range("K85").Value = 1059.25
$ws.Range("L85").Value = 1861
$ws.Range("M85").Value = 188.75
$ws.Range("N85").Value = -4357
$ws.Range("H113").Value = 2827.6365
$ws.Range("I113").Value = 2806.125
$ws.Range("K113").Value = 2806.125
$ws.Range("M113").Value = -636.125
$ws.Range("H132").Value = 337199.66
$ws.Range("I132").Value = 0
$ws.Range("K132").Value = 0
$ws.Range("M132").ClearContents()
$ws.Range("H136").Value = 63908.707
$ws.Range("J136").Value = 205560
$ws.Range("L136").Value = 616680
$ws.Range("N136").Value = -621780

# --- WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 2726.25
$ws.Range("I122").Value = 0
$ws.Range("J122").Value = 2726.25
$ws.Range("K122").Value = 0
$ws.Range("L122").Value = 8178.75
$ws.Range("M122").ClearContents()
$ws.Range("N122").Value = -13078.75
$ws.Range("H132").Value = 48922.383
$ws.Range("I132").Value = 29577.457
$ws.Range("J132").Value = 145647
$ws.Range("K132").Value = 88732.371
$ws.Range("L132").Value = 436941
$ws.Range("M132").Value = -86202.371
$ws.Range("N132").Value = -442001
